{"js": "// Replace the \"dividend\u00f7divisor=\" expressions in the worksheet table with\n// a new set of equivalent division problems (regenerated output).\nconst replacements = [\n  [\"336\u00f75=\", \"877\u00f74=\"],\n  [\"821\u00f77=\", \"263\u00f79=\"],\n  [\"675\u00f77=\", \"296\u00f79=\"],\n  [\"149\u00f73=\", \"357\u00f74=\"],\n  [\"935\u00f75=\", \"379\u00f78=\"],\n  [\"885\u00f78=\", \"645\u00f77=\"],\n  [\"624\u00f75=\", \"627\u00f74=\"],\n  [\"658\u00f79=\", \"838\u00f74=\"],\n  [\"481\u00f73=\", \"102\u00f74=\"],\n  [\"742\u00f77=\", \"998\u00f79=\"],\n  [\"446\u00f78=\", \"417\u00f76=\"],\n  [\"843\u00f73=\", \"915\u00f72=\"],\n  [\"100\u00f74=\", \"170\u00f79=\"],\n  [\"438\u00f78=\", \"177\u00f79=\"],\n  [\"297\u00f75=\", \"556\u00f76=\"],\n  [\"641\u00f74=\", \"826\u00f74=\"],\n  [\"320\u00f72=\", \"583\u00f77=\"],\n  [\"610\u00f79=\", \"219\u00f75=\"],\n  [\"122\u00f75=\", \"725\u00f76=\"],\n  [\"153\u00f78=\", \"535\u00f73=\"],\n  [\"193\u00f78=\", \"744\u00f74=\"],\n  [\"343\u00f79=\", \"213\u00f73=\"],\n  [\"407\u00f77=\", \"370\u00f75=\"],\n  [\"381\u00f79=\", \"896\u00f75=\"],\n  [\"143\u00f74=\", \"383\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the \"dividend\u00f7divisor=\" expressions in the worksheet table with a\n# new set of equivalent division problems (regenerated output).\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"336\u00f75=\", \"877\u00f74=\"),\n    @(\"821\u00f77=\", \"263\u00f79=\"),\n    @(\"675\u00f77=\", \"296\u00f79=\"),\n    @(\"149\u00f73=\", \"357\u00f74=\"),\n    @(\"935\u00f75=\", \"379\u00f78=\"),\n    @(\"885\u00f78=\", \"645\u00f77=\"),\n    @(\"624\u00f75=\", \"627\u00f74=\"),\n    @(\"658\u00f79=\", \"838\u00f74=\"),\n    @(\"481\u00f73=\", \"102\u00f74=\"),\n    @(\"742\u00f77=\", \"998\u00f79=\"),\n    @(\"446\u00f78=\", \"417\u00f76=\"),\n    @(\"843\u00f73=\", \"915\u00f72=\"),\n    @(\"100\u00f74=\", \"170\u00f79=\"),\n    @(\"438\u00f78=\", \"177\u00f79=\"),\n    @(\"297\u00f75=\", \"556\u00f76=\"),\n    @(\"641\u00f74=\", \"826\u00f74=\"),\n    @(\"320\u00f72=\", \"583\u00f77=\"),\n    @(\"610\u00f79=\", \"219\u00f75=\"),\n    @(\"122\u00f75=\", \"725\u00f76=\"),\n    @(\"153\u00f78=\", \"535\u00f73=\"),\n    @(\"193\u00f78=\", \"744\u00f74=\"),\n    @(\"343\u00f79=\", \"213\u00f73=\"),\n    @(\"407\u00f77=\", \"370\u00f75=\"),\n    @(\"381\u00f79=\", \"896\u00f75=\"),\n    @(\"143\u00f74=\", \"383\u00f74=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
